# Reorder the "proponente" entry in the columns list (column B) so that it
# moves from the top of the list (right after "descricao") down to the end
# of the list (right before "prop_pk1" / "prop_pk2"). This is done by
# shifting the values in B2:B15 up by one row and placing "proponente" in
# the now-vacant last slot (B15), effectively moving "proponente" to the
# bottom so it can be matched/aligned with the other proponentes automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    "processo",
    "emenda",
    "valor",
    "pontos_livre_(18_meses",
    "pontos_gesac",
    "pontos_indicados",
    "pontos_analisados",
    "pontos_aprovados",
    "encaminhamento",
    "data_aspar_informada",
    "data_cadastrado",
    "responsavel",
    "pendencia_28/12",
    "proponente"
)

$startRow = 2
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("B$row").Value = $newValues[$i]
}
